# Sprint #1 - Sprint Backlog: burndown chart update / prototype task removal
#
# The second burndown table (rows 17-24, the "Actual" chart) had a
# "Dimitar Stratiev" prototype task that was removed - its contribution
# to the K(=11th)/L(=12th)/M/N/O day columns collapses to 0 for the
# affected rows, which rolls up into the Daily Work Sum row (27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (Spike #1: Learn Zotero System) - Day 11-14 work removed
$ws.Range("M18:P18").Value = 0

# Row 19 (Spike #2: Learn JavaScript) - Day 11 work removed
$ws.Range("M19").Value = 0

# Row 21 (Task #2: create a checklist in the pop-up prompt) - Day 11-13 work removed
$ws.Range("M21:O21").Value = 0

# Row 23 (Task #4: implement file I/O) - Day 11-13 work removed
$ws.Range("M23:O23").Value = 0

# Row 24 (Task #5: documentation of functions) - Day 13 work removed
$ws.Range("O24").Value = 0

# Row 27 (Daily Work Sum) - manually tallied totals for Day 11 and Day 14
# (columns M and P aren't formulas in this sheet, so they're corrected by hand
#  to match the rest of the row, which recalculates automatically)
$ws.Range("M27").Value = 2.5
$ws.Range("P27").Value = 5.5

# Leave the selection where the editor finished working
$ws.Range("M28").Select()
